$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 15 (EmpNoDesc / 經辦合理性說明): length 50 -> 150, add remark note
$ws.Range("E15").Value = 150
$ws.Range("G15").Value = "2022/8/25長度放大150"

# Row 19 (ManagerDesc / 主管覆核說明): length 50 -> 150, add remark note, adjust formatting
$ws.Range("E19").Value = 150

$ws.Range("F15").Copy()
$ws.Range("F19").PasteSpecial(-4122)

$ws.Range("C17").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("G19").Value = "2022/8/25長度放大150"

$excel.CutCopyMode = $false

# Activate the DBD sheet and move the selection to B24, matching the saved view state
$ws.Activate()
$ws.Range("B24").Select()
